# Updates odds values in the betting sheet to reflect the latest quotes
# (Atualizando o arquivo XLSX)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.56
$ws.Range("H2").Value = 5.5
$ws.Range("I2").Value = 6.8
$ws.Range("J2").Value = 4.3
$ws.Range("K2").Value = 5.1
$ws.Range("N2").Value = 4
$ws.Range("O2").Value = 1.27
$ws.Range("Q2").Value = 1.79
$ws.Range("R2").Value = 1.41
$ws.Range("S2").Value = 3
$ws.Range("T2").Value = 1.87
$ws.Range("U2").Value = 1.94
$ws.Range("AF2").Value = 12
$ws.Range("AG2").Value = 12.5
$ws.Range("AN2").Value = 10.5
$ws.Range("F3").Value = 7.8
$ws.Range("G3").Value = 9.199999999999999
$ws.Range("I3").Value = 1.46
$ws.Range("J3").Value = 5.1
$ws.Range("L3").Value = 1.24
$ws.Range("M3").Value = 1.03
$ws.Range("P3").Value = 2.42
$ws.Range("R3").Value = 1.56
$ws.Range("S3").Value = 2.46
$ws.Range("T3").Value = 1.86
$ws.Range("U3").Value = 1.97
$ws.Range("V3").Value = 3.15
$ws.Range("X3").Value = 32
$ws.Range("Y3").Value = 12.5
$ws.Range("Z3").Value = 11.5
$ws.Range("AA3").Value = 13
$ws.Range("AB3").Value = 38
$ws.Range("AC3").Value = 15.5
$ws.Range("AD3").Value = 12.5
$ws.Range("AE3").Value = 15
$ws.Range("AG3").Value = 38
$ws.Range("AH3").Value = 26
$ws.Range("AI3").Value = 38
$ws.Range("AO3").Value = 6.6
$ws.Range("G4").Value = 26
$ws.Range("H4").Value = 1.12
$ws.Range("J4").Value = 1.32
$ws.Range("W4").Value = 1.04
$ws.Range("N5").Value = 2.66
$ws.Range("R5").Value = 1.18
$ws.Range("S5").Value = 2.68
$ws.Range("T5").Value = 1.04
$ws.Range("U5").Value = 1.04
$ws.Range("F6").Value = 1.76
$ws.Range("G6").Value = 1.86
$ws.Range("H6").Value = 5.1
$ws.Range("I6").Value = 6.2
$ws.Range("J6").Value = 3.35
$ws.Range("M6").Value = 1.07
$ws.Range("N6").Value = 2.98
$ws.Range("P6").Value = 1.68
$ws.Range("Q6").Value = 2.18
$ws.Range("R6").Value = 1.25
$ws.Range("S6").Value = 4.1
$ws.Range("T6").Value = 2.02
$ws.Range("U6").Value = 1.8
$ws.Range("V6").Value = 1.19
$ws.Range("W6").Value = 2.16
$ws.Range("X6").Value = 14
$ws.Range("Y6").Value = 18.5
$ws.Range("AB6").Value = 8.6
$ws.Range("AC6").Value = 9.800000000000001
$ws.Range("AD6").Value = 27
$ws.Range("AF6").Value = 12
$ws.Range("AG6").Value = 12.5
$ws.Range("AN6").Value = 19.5
